$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$nm = $p.NotesMaster
$tcsM = $m.Theme.ThemeColorScheme
$tcsN = $nm.Theme.ThemeColorScheme
for ($i=1; $i -le 12; $i++) {
    Write-Output "i=$i master=$($tcsM.Item($i).RGB) notes=$($tcsN.Item($i).RGB)"
}
